$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting rows 115:168 down to 116:169
$ws.Rows.Item(115).Insert()

# Populate the new row 115 with data
$ws.Cells.Item(115, 1).Value = 8
$ws.Cells.Item(115, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(115, 3).Value = "Coquimbo"
$ws.Cells.Item(115, 4).Value = 44510
$ws.Cells.Item(115, 5).Value = 4
$ws.Cells.Item(115, 6).Value = 100112003
$ws.Cells.Item(115, 7).Value = "Ajo"
$ws.Cells.Item(115, 8).Value = "Chino"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 600
$ws.Cells.Item(115, 11).Value = 19000
$ws.Cells.Item(115, 12).Value = 20000
$ws.Cells.Item(115, 13).Value = 19500
$ws.Cells.Item(115, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(115, 15).Value = "China"
$ws.Cells.Item(115, 16).Value = 1950
$ws.Cells.Item(115, 17).Value = 10
$ws.Cells.Item(115, 18).Value = "Hortaliza"

# Copy the style of column D (date format) from row 116 to row 115
$ws.Cells.Item(116, 4).Copy()
$ws.Cells.Item(115, 4).PasteSpecial(-4122) | Out-Null
